# Applies the price/volume refresh for cryptos.xlsx described in the commit
# "Updated cryptos list ... with GitHub Actions". Row 37/38 additionally swap
# (EthereumClassic <-> Monero) along with their refreshed price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.283.25"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.647.27"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'157.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "3.130.20"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "68.233.78"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "2.659.44"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D19").Value = "'364.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'75.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'558.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").Value = "'8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'1.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'160.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("E42").Value = "  +4.73%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'159.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'22.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "'0.614"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("E51").Value = "  +0.88%  "
